# Natmi following Dr Hou advice
#
# The sending/target cluster set grows from {FAPs, sCs} to {FAPs, sCs, ECs},
# so the 2x2 cluster grid (rows 2-5) becomes a 3x3 grid (rows 2-10) and all
# of the NATMI edge-weight / specificity statistics are recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T for data rows 2-10 (header row 1 is untouched).
# A = Sending cluster, B = Ligand symbol, C = Receptor symbol, D = Target cluster,
# E:T = the NATMI numeric statistics.
$rows = @(
    @("ECs",  "Crlf1", "Cntfr", "ECs",  2, 0.6666666666666666, 0.3448526666666667, 1.034558,           0.02513720817446796, 0.02513720817446796, 1, 0.3333333333333333, 0.05473433333333333, 0.164203,   0.005873455178951466, 0.005873455178951466, 0.01887528080822222, 0.169877527274,     0.00014764226553671,  0.00014764226553671),
    @("ECs",  "Crlf1", "Cntfr", "FAPs", 2, 0.6666666666666666, 0.3448526666666667, 1.034558,           0.02513720817446796, 0.02513720817446796, 3, 1,                  8.501068666666667,   25.503206,  0.912236301167251,    0.912236301167251,    2.931616199216445,   26.384545792948,    0.02293107380674783,  0.02293107380674783),
    @("ECs",  "Crlf1", "Cntfr", "sCs",  2, 0.6666666666666666, 0.3448526666666667, 1.034558,           0.02513720817446796, 0.02513720817446796, 3, 1,                  0.7631296666666666,  2.289389,   0.08189024365379753,  0.08189024365379755,  0.2631673005624444,  2.368505705062,     0.002058492102183412, 0.002058492102183413),
    @("FAPs", "Crlf1", "Cntfr", "ECs",  3, 1,                  8.096090333333334, 24.288271,           0.5901450902944959,  0.5901450902944958,  1, 0.3333333333333333, 0.05473433333333333, 0.164203,   0.005873455178951466, 0.005873455178951466, 0.4431341070014445,  3.988206963013,     0.003466190736922988, 0.003466190736922987),
    @("FAPs", "Crlf1", "Cntfr", "FAPs", 3, 1,                  8.096090333333334, 24.288271,           0.5901450902944959,  0.5901450902944958,  3, 1,                  8.501068666666667,   25.503206,  0.912236301167251,    0.912236301167251,    68.8254198552029,    619.428778696826,   0.5383517743222643,  0.5383517743222642),
    @("FAPs", "Crlf1", "Cntfr", "sCs",  3, 1,                  8.096090333333334, 24.288271,           0.5901450902944959,  0.5901450902944958,  3, 1,                  0.7631296666666666,  2.289389,   0.08189024365379753,  0.08189024365379755,  6.178366717379889,   55.605300456419,    0.04832712523530862, 0.04832712523530861),
    @("sCs",  "Crlf1", "Cntfr", "ECs",  3, 1,                  5.277870333333333, 15.833611,           0.3847177015310362,  0.3847177015310362,  1, 0.3333333333333333, 0.05473433333333333, 0.164203,   0.005873455178951466, 0.005873455178951466, 0.2888807141147778,  2.599926427033,     0.002259622176491769, 0.002259622176491769),
    @("sCs",  "Crlf1", "Cntfr", "FAPs", 3, 1,                  5.277870333333333, 15.833611,           0.3847177015310362,  0.3847177015310362,  3, 1,                  8.501068666666667,   25.503206,  0.912236301167251,    0.912236301167251,    44.86753811742955,   403.807843056866,   0.3509534530382389,  0.3509534530382389),
    @("sCs",  "Crlf1", "Cntfr", "sCs",  3, 1,                  5.277870333333333, 15.833611,           0.3847177015310362,  0.3847177015310362,  3, 1,                  0.7631296666666666,  2.289389,   0.08189024365379753,  0.08189024365379755,  4.027699428186556,   36.249294853679,    0.03150462631630551, 0.03150462631630552)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $rows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
